$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24 (table body row). This shifts the
# existing rows 24 ("Assignment Q&A Week 5"), 25 ("Labs 4") and 26
# ("ADA Compliance") down to 25, 26, 27 respectively while keeping each
# row's own formatting attached to its content (matches native Excel
# "Insert Table Rows Above" behaviour).
$ws.Rows("24:24").Insert()

# Grow the table (ListObject) to cover the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B2:F27"))

# Rename the existing "MVC - Core Sports Store App, 1" assignment (now on
# row 23) to reference Chapter 7 specifically.
$ws.Range("B23").Value = "MVC - Core Sports Store App, 1 - Chapter 7"

# Populate the freshly inserted row 24 with the new Chapter 8 entry.
$ws.Range("B24").Value = "MVC - Core Sports Store App, 1 - Chapter 8"
$ws.Range("C24").Value = 8
# Re-use the exact due-date text (incl. its leading non-breaking space)
# already stored in the shared-strings table instead of retyping it.
$ws.Range("D24").Value = $ws.Range("D23").Value2

# The "Assignment Q&A Week 5" row (shifted down to row 25) now records 1
# hour spent.
$ws.Range("C25").Value = 1

# Update the view state: scrolled so row 7 is at the top, with C21 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("C21").Select()
